# Atualização automática de TENENTE_PORTELA.xlsx
#
# - Rename "Paineis DARQ" -> "PAINEIS DARQ"
# - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Remove the "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

# Avoid any confirmation prompt when deleting a worksheet.
$excel.DisplayAlerts = $false | Out-Null

# Uppercase the two sheet names that changed casing/wording.
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Drop the now-unused "Desarquivamentos Pendentes" sheet.
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null

$excel.DisplayAlerts = $true | Out-Null
